# fix: fixed formatting when scrapping floating point numbers
#
# The "Importe" column (H2:H5) was scraped with Spanish/European-style
# number formatting (thousands separator "." and decimal separator ","),
# e.g. "2.710,00". This corrects it to use plain dot-decimal notation,
# e.g. "2710.00". The values remain plain text (as originally scraped),
# they are not converted into real numeric cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force the cells to Text format so Excel does not
# auto-convert the dot-decimal strings into numeric values when we
# assign them.
$ws.Range("H2:H5").NumberFormat = "@"

$ws.Range("H2").Value = "2710.00"
$ws.Range("H3").Value = "23576.88"
$ws.Range("H4").Value = "31435.84"
$ws.Range("H5").Value = "316.00"

# Restore the cells to the workbook's default (unformatted) style now
# that the text values are safely stored, so no visible formatting
# change is introduced.
$ws.Range("H2:H5").Style = "Normal"
